# 3DES Projetos - Adicionada aula08 - Desavio
#
# The PONTOS sheet gets a new column of grades (E, header "AT1-Cont")
# for a new class date (aula08, 2021-01-25 == serial 44221), and becomes
# the active/selected sheet (instead of FREQ).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PONTOS")

# New class date header value for column E, row 2 (keeps existing date style).
$ws.Range("E2").Value = 44221

# New column of per-student scores for the aula08 assessment.
$scores = @{
    3  = 0
    4  = 0
    5  = 0
    6  = 2
    7  = 0
    8  = 2
    9  = 2
    10 = 3
    11 = 0
    12 = 2
    13 = 0
    14 = 3
    15 = 4
    16 = 2
    17 = 2
    18 = 2
    19 = 3
    20 = 2
}

foreach ($row in $scores.Keys) {
    $ws.Cells.Item($row, 5).Value = $scores[$row]
}

# Make PONTOS the active sheet/tab (was FREQ before), and move its
# selection to E21.
$ws.Activate()
$ws.Range("E21").Select()
